$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$msoTextOrientationHorizontal = 1
$msoAutoSizeShapeToFitText    = 1
$ppAlignCenter                = 2
$EMU_PER_POINT                = 12700

# Target position/size (from the OOXML diff), expressed in points since the
# PowerPoint object model's Shape geometry properties are in points.
$left   = 29210184 / $EMU_PER_POINT
$top    = 17118396 / $EMU_PER_POINT
$width  = 11761378 / $EMU_PER_POINT
$height = 1415772  / $EMU_PER_POINT

# Add the new "UML Diagram" caption text box below/near the UML picture.
$tb = $s.Shapes.AddTextbox($msoTextOrientationHorizontal, $left, $top, $width, $height)
$tb.Name = "TextBox 9"
$tb.Fill.Visible = $false

$tf = $tb.TextFrame
$tf.WordWrap = $true
$tf.AutoSize = $msoAutoSizeShapeToFitText

$tr = $tf.TextRange
$tr.Text = "UML Diagram"
$tr.ParagraphFormat.Alignment = $ppAlignCenter

# Re-assert the exact geometry: setting the text can trigger shape-to-fit-text
# autofit, which recomputes the height from the rendered line(s) of text.
$tb.Left   = $left
$tb.Top    = $top
$tb.Width  = $width
$tb.Height = $height
